$d = $word.ActiveDocument

# --- Programa (PT) ---
$f = $d.Content
$found = $f.Find.Execute("Conceitos introdutórios: adsorvente e adsorvato. Adsorventes orgânicos e inorgânicos. Características dos adsorventes: sítios de adsorção, área superficial, porosidadeAdsorvatos catiônicos e aniônicos. Condicionantes do processo de adsorção.Cinética da reação de adsorçãoMecanismos de adsorção. Análise de isotermas de Langmuir, Freundlich, Temkin e SipsExemplos de ocorrências e aplicações do processo de adsorção em Engenharia AmbientalAulas práticas", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text for Programa (PT)"
}
$target = $d.Range($f.Start, $f.End)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Conceitos introdutórios: adsorvente e adsorvato. </w:t><w:br/><w:t>Adsorventes orgânicos e inorgânicos. Características dos adsorventes: sítios de adsorção, área superficial, porosidade</w:t><w:br/><w:t xml:space="preserve">Adsorvatos catiônicos e aniônicos. </w:t><w:br/><w:t>Condicionantes do processo de adsorção.</w:t><w:br/><w:t>Cinética da reação de adsorção</w:t><w:br/><w:t>Mecanismos de adsorção. Análise de isotermas de Langmuir, Freundlich, Temkin e Sips</w:t><w:br/><w:t>Exemplos de ocorrências e aplicações do processo de adsorção em Engenharia Ambiental</w:t><w:br/><w:t>Aulas práticas</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml)

# --- Programa (EN) ---
$f = $d.Content
$found = $f.Find.Execute("Introductory concepts: adsorbent and adsorbate. Organic and inorganic adsorbentsCharacteristics of adsorbents: adsorption sites, surface area, porosity Cationic and anionic adsorbates. Conditions of the adsorption process. Kinetics of the adsorption reaction Adsorption mechanisms. Analysis of Langmuir, Freundlich, Temkin and Sips isotherms Examples of occurrences and applications of the adsorption process in Environmental Engineering Practical classes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text for Programa (EN)"
}
$target = $d.Range($f.Start, $f.End)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Introductory concepts: adsorbent and adsorbate. </w:t><w:br/><w:t>Organic and inorganic adsorbents</w:t><w:br/><w:t xml:space="preserve">Characteristics of adsorbents: adsorption sites, surface area, porosity </w:t><w:br/><w:t xml:space="preserve">Cationic and anionic adsorbates. </w:t><w:br/><w:t xml:space="preserve">Conditions of the adsorption process. </w:t><w:br/><w:t xml:space="preserve">Kinetics of the adsorption reaction </w:t><w:br/><w:t xml:space="preserve">Adsorption mechanisms. Analysis of Langmuir, Freundlich, Temkin and Sips isotherms </w:t><w:br/><w:t>Examples of occurrences and applications of the adsorption process in Environmental Engineering Practical classes</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml)

# --- Bibliografia ---
$f = $d.Content
$found = $f.Find.Execute("BOSCOV, M. E. Geotecnia ambiental. Oficina de Textos: São Paulo, 2008. 248p.YONG, R. N. Geoenvironmental engineering. contaminated soils, Pollutant fate and migration. CRC Press, 2001. 307p.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text for Bibliografia"
}
$target = $d.Range($f.Start, $f.End)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>BOSCOV, M. E. Geotecnia ambiental. Oficina de Textos: São Paulo, 2008. 248p.</w:t><w:br/><w:t>YONG, R. N. Geoenvironmental engineering. contaminated soils, Pollutant fate and migration. CRC Press, 2001. 307p.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml)

Write-Host "All replacements complete"
